$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aircraft_scheduling")

# "Clean up" the E column: replace the formulas (=Cx+offset) that overflowed
# past midnight with the plain decimal time-of-day value that was intended
# (i.e. drop the extra day). This matches what F (=MOD(E-C,1)) already
# expects, so F recalculates automatically once E no longer carries a
# formula.
$ws.Range("E4").Value  = 0.24111111111111114
$ws.Range("E5").Value  = 0.16888888888888889
$ws.Range("E9").Value  = 0.21902777777777779
$ws.Range("E11").Value = 0.4319675925925926
$ws.Range("E12").Value = 0.34127314814814813
$ws.Range("E13").Value = 0.52976851851851847
$ws.Range("E14").Value = 0.70796296296296291
$ws.Range("E15").Value = 0.51861111111111113
$ws.Range("E16").Value = 0.62513888888888891
$ws.Range("E17").Value = 0.61722222222222223
$ws.Range("E18").Value = 0.77541666666666664

# Reflect the author's last cursor position on this sheet.
$ws.Activate()
$ws.Range("E19").Select()
